$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Rename the classifier header "Regras" -> "Odds" (E5 / N5 mirror header, and bottom row 37)
$ws.Range("E5").Value = "Odds"
$ws.Range("N5").Value = "Odds"
$ws.Range("E37").Value = "Odds"

# Fill in previously empty "Odds" score column (E6:E35) with the new values
$ws.Range("E6").Value = 0.55419135587312596
$ws.Range("E7").Value = 0.55415214360404297
$ws.Range("E8").Value = 0.55404757755315404
$ws.Range("E9").Value = 0.55420442662948699
$ws.Range("E10").Value = 0.55421314046706105
$ws.Range("E11").Value = 0.55420878354827396
$ws.Range("E12").Value = 0.55418264203555201
$ws.Range("E13").Value = 0.55407371906587599
$ws.Range("E14").Value = 0.55416957127919098
$ws.Range("E15").Value = 0.554165214360404
$ws.Range("E16").Value = 0.55414342976646902
$ws.Range("E17").Value = 0.55415214360404297
$ws.Range("E18").Value = 0.55418264203555201
$ws.Range("E19").Value = 0.55415650052283005
$ws.Range("E20").Value = 0.55416085744161703
$ws.Range("E21").Value = 0.55419135587312596
$ws.Range("E22").Value = 0.554221854304635
$ws.Range("E23").Value = 0.55418699895433898
$ws.Range("E24").Value = 0.55418699895433898
$ws.Range("E25").Value = 0.55418699895433898
$ws.Range("E26").Value = 0.55411293133495998
$ws.Range("E27").Value = 0.55415214360404297
$ws.Range("E28").Value = 0.55416085744161703
$ws.Range("E29").Value = 0.55407807598466297
$ws.Range("E30").Value = 0.55418699895433898
$ws.Range("E31").Value = 0.55410857441617201
$ws.Range("E32").Value = 0.55413907284768205
$ws.Range("E33").Value = 0.55414778668525599
$ws.Range("E34").Value = 0.55415214360404297
$ws.Range("E35").Value = 0.55414342976646902

$excel.CalculateFullRebuild()

# Update sheet view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("H10").Select()
